$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 36.40130033333333
$ws.Range("H2").Value = 109.203901
$ws.Range("I2").Value = 0.1897437225523226
$ws.Range("J2").Value = 0.1897437225523226
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 99.15720277683388
$ws.Range("R2").Value = 892.414824991505
$ws.Range("S2").Value = 0.008795571158747078
$ws.Range("T2").Value = 0.00879557115874708
$ws.Range("G3").Value = 36.40130033333333
$ws.Range("H3").Value = 109.203901
$ws.Range("I3").Value = 0.1897437225523226
$ws.Range("J3").Value = 0.1897437225523226
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 1478.643776626744
$ws.Range("R3").Value = 13307.79398964069
$ws.Range("S3").Value = 0.1311605833116294
$ws.Range("T3").Value = 0.1311605833116294
$ws.Range("G4").Value = 36.40130033333333
$ws.Range("H4").Value = 109.203901
$ws.Range("I4").Value = 0.1897437225523226
$ws.Range("J4").Value = 0.1897437225523226
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 561.2820242102616
$ws.Range("R4").Value = 5051.538217892355
$ws.Range("S4").Value = 0.04978756808194616
$ws.Range("T4").Value = 0.04978756808194616
$ws.Range("I5").Value = 0.6107553255746098
$ws.Range("J5").Value = 0.6107553255746098
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 319.1715059154745
$ws.Range("R5").Value = 2872.54355323927
$ws.Range("S5").Value = 0.02831156601343628
$ws.Range("T5").Value = 0.02831156601343628
$ws.Range("I6").Value = 0.6107553255746098
$ws.Range("J6").Value = 0.6107553255746098
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("R6").Value = 42835.70460984058
$ws.Range("S6").Value = 0.4221853755449542
$ws.Range("T6").Value = 0.4221853755449542
$ws.Range("I7").Value = 0.6107553255746098
$ws.Range("J7").Value = 0.6107553255746098
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 1806.678928949464
$ws.Range("R7").Value = 16260.11036054517
$ws.Range("S7").Value = 0.1602583840162193
$ws.Range("T7").Value = 0.1602583840162193
$ws.Range("G8").Value = 38.27317166666666
$ws.Range("H8").Value = 114.819515
$ws.Range("I8").Value = 0.1995009518730676
$ws.Range("J8").Value = 0.1995009518730676
$ws.Range("M8").Value = 2.724001666666667
$ws.Range("N8").Value = 8.172005
$ws.Range("O8").Value = 0.04635500474236593
$ws.Range("P8").Value = 0.04635500474236593
$ws.Range("Q8").Value = 104.2561834086194
$ws.Range("R8").Value = 938.305650677575
$ws.Range("S8").Value = 0.009247867570182568
$ws.Range("T8").Value = 0.009247867570182568
$ws.Range("G9").Value = 38.27317166666666
$ws.Range("H9").Value = 114.819515
$ws.Range("I9").Value = 0.1995009518730676
$ws.Range("J9").Value = 0.1995009518730676
$ws.Range("O9").Value = 0.6912512390256352
$ws.Range("P9").Value = 0.6912512390256351
$ws.Range("Q9").Value = 1554.680370713598
$ws.Range("R9").Value = 13992.12333642238
$ws.Range("S9").Value = 0.1379052801690516
$ws.Range("T9").Value = 0.1379052801690516
$ws.Range("G10").Value = 38.27317166666666
$ws.Range("H10").Value = 114.819515
$ws.Range("I10").Value = 0.1995009518730676
$ws.Range("J10").Value = 0.1995009518730676
$ws.Range("M10").Value = 15.419285
$ws.Range("N10").Value = 46.257855
$ws.Range("O10").Value = 0.2623937562319988
$ws.Range("P10").Value = 0.2623937562319988
$ws.Range("Q10").Value = 590.1449417822583
$ws.Range("R10").Value = 5311.304476040325
$ws.Range("S10").Value = 0.05234780413383344
$ws.Range("T10").Value = 0.05234780413383344
